# update scripts wuth new tpm
# Recomputed NATMI ligand-receptor (Vcam1-Itga4) TPM-derived stats for rows 2-10 (columns G:T).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.775841999999999
$ws.Range("H2").Value = 8.327525999999999
$ws.Range("I2").Value = 0.0624750527258915
$ws.Range("J2").Value = 0.0624750527258915
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.023286
$ws.Range("N2").Value = 0.069858
$ws.Range("O2").Value = 0.009310710475795457
$ws.Range("P2").Value = 0.009310710475795458
$ws.Range("Q2").Value = 0.06463825681199999
$ws.Range("R2").Value = 0.581744311308
$ws.Range("S2").Value = 0.0005816871278908315
$ws.Range("T2").Value = 0.0005816871278908316

# Row 3
$ws.Range("G3").Value = 2.775841999999999
$ws.Range("H3").Value = 8.327525999999999
$ws.Range("I3").Value = 0.0624750527258915
$ws.Range("J3").Value = 0.0624750527258915
$ws.Range("O3").Value = 0.05314667307834813
$ws.Range("P3").Value = 0.05314667307834814
$ws.Range("Q3").Value = 0.3689630680786665
$ws.Range("R3").Value = 3.320667612707999
$ws.Range("S3").Value = 0.003320341202775518
$ws.Range("T3").Value = 0.003320341202775518

# Row 4
$ws.Range("G4").Value = 2.775841999999999
$ws.Range("H4").Value = 8.327525999999999
$ws.Range("I4").Value = 0.0624750527258915
$ws.Range("J4").Value = 0.0624750527258915
$ws.Range("M4").Value = 2.344785333333334
$ws.Range("N4").Value = 7.034356000000001
$ws.Range("O4").Value = 0.9375426164458565
$ws.Range("P4").Value = 0.9375426164458565
$ws.Range("Q4").Value = 6.508753609250666
$ws.Range("R4").Value = 58.578782483256
$ws.Range("S4").Value = 0.05857302439522515
$ws.Range("T4").Value = 0.05857302439522516

# Row 5
$ws.Range("I5").Value = 0.2652892219050753
$ws.Range("J5").Value = 0.2652892219050753
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.023286
$ws.Range("N5").Value = 0.069858
$ws.Range("O5").Value = 0.009310710475795457
$ws.Range("P5").Value = 0.009310710475795458
$ws.Range("Q5").Value = 0.274474884082
$ws.Range("R5").Value = 2.470273956738
$ws.Range("S5").Value = 0.00247003113750721
$ws.Range("T5").Value = 0.00247003113750721

# Row 6
$ws.Range("I6").Value = 0.2652892219050753
$ws.Range("J6").Value = 0.2652892219050753
$ws.Range("O6").Value = 0.05314667307834813
$ws.Range("P6").Value = 0.05314667307834814
$ws.Range("S6").Value = 0.01409923954779839
$ws.Range("T6").Value = 0.01409923954779839

# Row 7
$ws.Range("I7").Value = 0.2652892219050753
$ws.Range("J7").Value = 0.2652892219050753
$ws.Range("M7").Value = 2.344785333333334
$ws.Range("N7").Value = 7.034356000000001
$ws.Range("O7").Value = 0.9375426164458565
$ws.Range("P7").Value = 0.9375426164458565
$ws.Range("Q7").Value = 27.63826687983512
$ws.Range("R7").Value = 248.744401918516
$ws.Range("S7").Value = 0.2487199512197697
$ws.Range("T7").Value = 0.2487199512197697

# Row 8
$ws.Range("G8").Value = 29.86824466666667
$ws.Range("H8").Value = 89.60473400000001
$ws.Range("I8").Value = 0.6722357253690333
$ws.Range("J8").Value = 0.6722357253690333
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.023286
$ws.Range("N8").Value = 0.069858
$ws.Range("O8").Value = 0.009310710475795457
$ws.Range("P8").Value = 0.009310710475795458
$ws.Range("Q8").Value = 0.6955119453080001
$ws.Range("R8").Value = 6.259607507772001
$ws.Range("S8").Value = 0.006258992210397416
$ws.Range("T8").Value = 0.006258992210397417

# Row 9
$ws.Range("G9").Value = 29.86824466666667
$ws.Range("H9").Value = 89.60473400000001
$ws.Range("I9").Value = 0.6722357253690333
$ws.Range("J9").Value = 0.6722357253690333
$ws.Range("O9").Value = 0.05314667307834813
$ws.Range("P9").Value = 0.05314667307834814
$ws.Range("Q9").Value = 3.970067168930222
$ws.Range("R9").Value = 35.730604520372
$ws.Range("S9").Value = 0.03572709232777423
$ws.Range("T9").Value = 0.03572709232777423

# Row 10
$ws.Range("G10").Value = 29.86824466666667
$ws.Range("H10").Value = 89.60473400000001
$ws.Range("I10").Value = 0.6722357253690333
$ws.Range("J10").Value = 0.6722357253690333
$ws.Range("M10").Value = 2.344785333333334
$ws.Range("N10").Value = 7.034356000000001
$ws.Range("O10").Value = 0.9375426164458565
$ws.Range("P10").Value = 0.9375426164458565
$ws.Range("Q10").Value = 70.03462202681158
$ws.Range("R10").Value = 630.3115982413041
$ws.Range("S10").Value = 0.6302496408308617
$ws.Range("T10").Value = 0.6302496408308617
